$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new trailing columns: L = "Village", M = "Agronomist" ---
# Copy the formatting of column K (rows 1-9, which already has the header
# style in rows 1-2 and the blank-body style in rows 3-9) into the two new
# columns so the new cells pick up matching borders/fill/number format.
$ws.Range("K1:K9").Copy()
$ws.Range("L1:L9").PasteSpecial(-4122)
$ws.Range("K1:K9").Copy()
$ws.Range("M1:M9").PasteSpecial(-4122)

# Header row (row 1)
$ws.Range("L1").Value = "Village"
$ws.Range("M1").Value = "Agronomist"

# Data row (row 2)
$ws.Range("L2").Value = "Ukpo"
$ws.Range("M2").Value = "Paul walker"

# Rows 3-9 stay blank (already formatted above) for columns L and M.

# --- Remove the now-unused trailing row 10 ---
$ws.Rows("10").Delete()
